$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.082.07"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.815.57"
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'310.52"
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").Value = "'0.4990"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("D8").Value = "'0.3905"
$ws.Range("E8").Value = "  -0.91%  "
$ws.Range("D9").Value = "'0.09875"
$ws.Range("E9").Value = "  +26.51%  "
$ws.Range("D10").Value = "'1.109"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'40.89"
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("D12").Value = "'6.424"
$ws.Range("E12").Value = "  +2.81%  "
$ws.Range("E13").Value = "  +1.64%  "
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "1.812.04"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("D16").Value = "'7.271"
$ws.Range("E16").Value = "  +0.41%  "
$ws.Range("D17").Value = "'0.00001140"
$ws.Range("E17").Value = "  +5.84%  "
$ws.Range("D18").Value = "'92.37"
$ws.Range("E18").Value = "  +0.84%  "
$ws.Range("D19").Value = "'0.06633"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("D22").Value = "'5.951"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "28.133.99"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").Value = "'2.241"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("D26").Value = "'159.19"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("D27").Value = "'20.69"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "2.022.69"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("D29").Value = "'2.405"
$ws.Range("E29").Value = "  +1.64%  "
$ws.Range("D30").Value = "'126.77"
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").Value = "'0.1058"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("D32").Value = "'1.035"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "'5.569"
$ws.Range("E33").Value = "  +1.34%  "
$ws.Range("D34").Value = "'3.621"
$ws.Range("E34").Value = "  -0.25%  "
$ws.Range("D35").Value = "'0.06680"
$ws.Range("E35").Value = "  -5.70%  "
$ws.Range("D36").Value = "'0.02344"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").Value = "'8.894"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("D38").Value = "'0.2142"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("D39").Value = "'4.957"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "'11.35"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "'0.6198"
$ws.Range("E41").Value = "  +1.40%  "
$ws.Range("D42").Value = "'1.180"
$ws.Range("E42").Value = "  +2.38%  "
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'13.15"
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("D45").Value = "'0.5901"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "'3.690"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'1.266"
$ws.Range("E47").Value = "  -3.47%  "
$ws.Range("D48").Value = "'124.31"
$ws.Range("E48").Value = "  -0.44%  "
$ws.Range("D49").Value = "'1.941"
$ws.Range("E50").Value = "  -2.28%  "
$ws.Range("D51").Value = "'0.06780"
$ws.Range("E51").Value = "  -0.52%  "
